$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1141.1538
$ws.Range("I41").Value = 461.7143
$ws.Range("J41").Value = 1933.8334
$ws.Range("K41").Value = 461.7143
$ws.Range("L41").Value = 1933.8334
$ws.Range("M41").Value = -21.71429999999998
$ws.Range("N41").Value = -2813.8334
$ws.Range("H43").Value = 18842128
$ws.Range("J43").Value = 1755.125
$ws.Range("L43").Value = 1755.125
$ws.Range("N43").Value = -1893.125
$ws.Range("H92").Value = 2232979.8
$ws.Range("I92").Value = 977400.0600000001
$ws.Range("J92").Value = 6250834.5
$ws.Range("K92").Value = 977400.0600000001
$ws.Range("L92").Value = 6250834.5
$ws.Range("M92").Value = -976152.0600000001
$ws.Range("N92").Value = -6253330.5
$ws.Range("H103").Value = 514.38464
$ws.Range("I103").Value = 361.25
$ws.Range("K103").Value = 1083.75
$ws.Range("M103").Value = -497.75
$ws.Range("H125").Value = 1682.5714
$ws.Range("I125").Value = 746.8333
$ws.Range("K125").Value = 6721.4997
$ws.Range("M125").Value = -4261.4997
$ws.Range("H132").Value = 30234.97
$ws.Range("I132").Value = 31999.688
$ws.Range("K132").Value = 95999.064
$ws.Range("M132").Value = -93469.064
$ws.Range("H133").Value = 147467.5
$ws.Range("J133").Value = 147467.5
$ws.Range("L133").Value = 147467.5
$ws.Range("N133").Value = -157587.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2214.2727
$ws.Range("I2").Value = 935.7
$ws.Range("K2").Value = 935.7
$ws.Range("M2").Value = -822.7
$ws.Range("H3").Value = 20005
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H32").Value = 220895.02
$ws.Range("I32").Value = 251039.67
$ws.Range("K32").Value = 251039.67
$ws.Range("M32").Value = -250752.67
$ws.Range("H116").Value = 2214.2727
$ws.Range("I116").Value = 935.7
$ws.Range("K116").Value = 935.7
$ws.Range("M116").Value = 1358.3
$ws.Range("H122").Value = 1332
$ws.Range("I122").Value = 1332
$ws.Range("K122").Value = 3996
$ws.Range("M122").Value = -1546
$ws.Range("H134").Value = 87319.25
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 87319.25
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 87319.25
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -97459.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2214.2727
$ws.Range("I3").Value = 935.7
$ws.Range("K3").Value = 935.7
$ws.Range("M3").Value = -821.7
$ws.Range("H7").Value = 18333700
$ws.Range("I7").Value = 27500000
$ws.Range("K7").Value = 27500000
$ws.Range("M7").Value = -27499887
$ws.Range("H86").Value = 3611.8572
$ws.Range("J86").Value = 8849.5
$ws.Range("L86").Value = 8849.5
$ws.Range("N86").Value = -11095.5
$ws.Range("H89").Value = 3611.8572
$ws.Range("J89").Value = 8849.5
$ws.Range("L89").Value = 44247.5
$ws.Range("N89").Value = -55479.5
$ws.Range("H107").Value = 6538.8047
$ws.Range("I107").Value = 7746.4062
$ws.Range("K107").Value = 7746.4062
$ws.Range("M107").Value = -5826.4062

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2669.7693
$ws.Range("I122").Value = 2791.1428
$ws.Range("K122").Value = 8373.428400000001
$ws.Range("M122").Value = -5923.428400000001
$ws.Range("H132").Value = 49458.617
$ws.Range("I132").Value = 49458.617
$ws.Range("K132").Value = 148375.851
$ws.Range("M132").Value = -145845.851
$ws.Range("H134").Value = 1580.0526
$ws.Range("I134").Value = 1241.6428
$ws.Range("K134").Value = 3724.9284
$ws.Range("M134").Value = -1189.9284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 14569.75
$ws.Range("J62").Value = 19640
$ws.Range("L62").Value = 58920
$ws.Range("N62").Value = -60292
$ws.Range("H63").Value = 7767.25
$ws.Range("J63").Value = 19166.666
$ws.Range("L63").Value = 57499.99800000001
$ws.Range("N63").Value = -58997.99800000001
$ws.Range("H65").Value = 14569.75
$ws.Range("J65").Value = 19640
$ws.Range("L65").Value = 176760
$ws.Range("N65").Value = -183624
$ws.Range("H66").Value = 7767.25
$ws.Range("J66").Value = 19166.666
$ws.Range("L66").Value = 172499.994
$ws.Range("N66").Value = -179987.994
$ws.Range("H69").Value = 9229.214
$ws.Range("I69").Value = 16305
$ws.Range("J69").Value = 6398.9
$ws.Range("K69").Value = 48915
$ws.Range("L69").Value = 19196.7
$ws.Range("M69").Value = -48104
$ws.Range("N69").Value = -20818.7
$ws.Range("H72").Value = 9229.214
$ws.Range("I72").Value = 16305
$ws.Range("J72").Value = 6398.9
$ws.Range("K72").Value = 146745
$ws.Range("L72").Value = 57590.1
$ws.Range("M72").Value = -142689
$ws.Range("N72").Value = -65702.10000000001
$ws.Range("H75").Value = 125001800
$ws.Range("I75").Value = 250001120
$ws.Range("J75").Value = 41668916
$ws.Range("K75").Value = 750003360
$ws.Range("L75").Value = 125006748
$ws.Range("M75").Value = -750002362
$ws.Range("N75").Value = -125008744
$ws.Range("H78").Value = 125001800
$ws.Range("I78").Value = 250001120
$ws.Range("J78").Value = 41668916
$ws.Range("K78").Value = 2250010080
$ws.Range("L78").Value = 375020244
$ws.Range("M78").Value = -2250005088
$ws.Range("N78").Value = -375030228
$ws.Range("H99").Value = 166677490
$ws.Range("I99").Value = 333334660
$ws.Range("K99").Value = 1000003980
$ws.Range("M99").Value = -1000001734
$ws.Range("H106").Value = 18323.334
$ws.Range("J106").Value = 18323.334
$ws.Range("L106").Value = 54970.00199999999
$ws.Range("N106").Value = -56862.00199999999
$ws.Range("H130").Value = 6681.3335
$ws.Range("J130").Value = 4999
$ws.Range("L130").Value = 14997
$ws.Range("N130").Value = -25037

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 64991
$ws.Range("J42").Value = 64991
$ws.Range("L42").Value = 64991
$ws.Range("N42").Value = -65961
$ws.Range("H70").Value = 7335.6665
$ws.Range("I70").Value = 4503.5
$ws.Range("J70").Value = 8751.75
$ws.Range("K70").Value = 4503.5
$ws.Range("L70").Value = 8751.75
$ws.Range("M70").Value = -4233.5
$ws.Range("N70").Value = -9291.75
$ws.Range("H73").Value = 7335.6665
$ws.Range("I73").Value = 4503.5
$ws.Range("J73").Value = 8751.75
$ws.Range("K73").Value = 4503.5
$ws.Range("L73").Value = 8751.75
$ws.Range("M73").Value = -3567.5
$ws.Range("N73").Value = -10623.75
$ws.Range("H80").Value = 15834.2
$ws.Range("I80").Value = 10365.6
$ws.Range("J80").Value = 24037.1
$ws.Range("K80").Value = 10365.6
$ws.Range("L80").Value = 24037.1
$ws.Range("M80").Value = -9367.6
$ws.Range("N80").Value = -26033.1
$ws.Range("H83").Value = 15834.2
$ws.Range("I83").Value = 10365.6
$ws.Range("J83").Value = 24037.1
$ws.Range("K83").Value = 51828
$ws.Range("L83").Value = 120185.5
$ws.Range("M83").Value = -46836
$ws.Range("N83").Value = -130169.5
$ws.Range("H102").Value = 13159226
$ws.Range("I102").Value = 16130327
$ws.Range("K102").Value = 16130327
$ws.Range("M102").Value = -16128705
$ws.Range("H115").Value = 64991
$ws.Range("J115").Value = 64991
$ws.Range("L115").Value = 64991
$ws.Range("N115").Value = -67341
$ws.Range("H122").Value = 2166
$ws.Range("I122").Value = 2206.4167
$ws.Range("K122").Value = 6619.250100000001
$ws.Range("M122").Value = -4169.250100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 16873.75
$ws.Range("J13").Value = 16873.75
$ws.Range("L13").Value = 16873.75
$ws.Range("N13").Value = -17153.75
$ws.Range("H16").Value = 1568.8948
$ws.Range("I16").Value = 1333.0667
$ws.Range("J16").Value = 2453.25
$ws.Range("K16").Value = 1333.0667
$ws.Range("L16").Value = 2453.25
$ws.Range("M16").Value = -1163.0667
$ws.Range("N16").Value = -2793.25
$ws.Range("H22").Value = 3947.6667
$ws.Range("J22").Value = 5666.9375
$ws.Range("L22").Value = 5666.9375
$ws.Range("N22").Value = -6256.9375
$ws.Range("H27").Value = 3947.6667
$ws.Range("J27").Value = 5666.9375
$ws.Range("L27").Value = 5666.9375
$ws.Range("N27").Value = -5880.9375
